# Kaman TC18_Verify_ShippingPage - "Changes done for Kaman new UI - header & footer"
#
# Summary of the edit:
#  Sheet "TC18_Verify_ShippingPage" (test-steps sheet):
#   - Row 14 Keyword SCROLL_DOWN            -> TINY_SCROLL_DOWN
#   - Row 19 (CLICK / ViewFullCart / CSS)   -> removed entirely (row deleted,
#     everything below shifts up by one row)
#   - The second SCROLL_DOWN row (old row 26, new row 25) -> TINY_SCROLL_DOWN
#
#  Sheet "Testdata" (object repository sheet):
#   - Two new rows appended at the bottom:
#       EleType1 | JSElement
#       EleType2 | JSElement

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---- Sheet1: TC18_Verify_ShippingPage ----

# Remove the "CLICK ViewFullCart" step (row 19) - everything below shifts up.
$ws1.Rows.Item(19).Delete()

# Rename both SCROLL_DOWN steps to TINY_SCROLL_DOWN (row 14 keeps its number,
# the second one used to be row 26 and is now row 25 after the deletion).
$ws1.Cells.Item(14, 2).Value = "TINY_SCROLL_DOWN"
$ws1.Cells.Item(25, 2).Value = "TINY_SCROLL_DOWN"

# ---- Sheet2: Testdata ----

# Append the two new object-repository entries.
$ws2.Cells.Item(35, 1).Value = "EleType1"
$ws2.Cells.Item(35, 2).Value = "JSElement"
$ws2.Cells.Item(36, 1).Value = "EleType2"
$ws2.Cells.Item(36, 2).Value = "JSElement"

# Match the bordered-cell look used by the rest of the table.
$ws2.Range("A35:B36").Borders.LineStyle = 1

# ---- Selections / active sheet, to mirror the saved view state ----

$ws2.Activate()
$ws2.Range("A35:B36").Select()

$ws1.Activate()
$ws1.Range("B25").Select()
